# Update "想去人数" (number of interested attendees) figures on the
# "展览" and "全部类型" sheets to reflect the latest scrape results.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 8336
    $ws.Range("F4").Value = 191
    $ws.Range("F5").Value = 347
}
